$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 158, shifting existing rows 158:235 down to 159:236.
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with the new observation.
$ws.Range("A158").Value = 10
$ws.Range("B158").Value = "Vega Modelo de Temuco"
$ws.Range("C158").Value = "La Araucanía"
$ws.Range("D158").Value = 44466
$ws.Range("E158").Value = 9
$ws.Range("F158").Value = 100112032
$ws.Range("G158").Value = "Zapallo italiano"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 100
$ws.Range("K158").Value = 17000
$ws.Range("L158").Value = 18000
$ws.Range("M158").Value = 17500
$ws.Range("N158").Value = "$/caja 60 unidades"
$ws.Range("O158").Value = "Región de Arica y Parinacota"
$ws.Range("P158").Value = 292
$ws.Range("Q158").Value = 60
$ws.Range("R158").Value = "Hortaliza"
